# Auto-generated edit script
# Applies numeric corrections to market-data columns (H-N) across all 8 job sheets
# as captured by the scheduled Aegis Profits runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5454.2856
$ws.Range("J69").Value = 5530
$ws.Range("L69").Value = 16590
$ws.Range("N69").Value = -18338
$ws.Range("H72").Value = 5454.2856
$ws.Range("J72").Value = 5530
$ws.Range("L72").Value = 49770
$ws.Range("N72").Value = -58506
$ws.Range("H132").Value = 10009213
$ws.Range("I132").Value = 10425800
$ws.Range("J132").Value = 11111
$ws.Range("K132").Value = 31277400
$ws.Range("L132").Value = 33333
$ws.Range("M132").Value = -31274870
$ws.Range("N132").Value = -38393
$ws.Range("H137").Value = 1313.8649
$ws.Range("I137").Value = 1134.6571
$ws.Range("K137").Value = 3403.9713
$ws.Range("M137").Value = -853.9712999999997

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20646.055
$ws.Range("I32").Value = 4400.2046
$ws.Range("K32").Value = 4400.2046
$ws.Range("M32").Value = -4113.2046
$ws.Range("H61").Value = 2121.372
$ws.Range("I61").Value = 1227.3158
$ws.Range("J61").Value = 2829.1667
$ws.Range("K61").Value = 1227.3158
$ws.Range("L61").Value = 2829.1667
$ws.Range("M61").Value = -1015.3158
$ws.Range("N61").Value = -3253.1667
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
$ws.Range("H113").Value = 36500
$ws.Range("J113").Value = 36500
$ws.Range("L113").Value = 36500
$ws.Range("N113").Value = -45178
$ws.Range("H124").Value = 28285.8
$ws.Range("J124").Value = 28285.8
$ws.Range("L124").Value = 28285.8
$ws.Range("N124").Value = -38105.8
$ws.Range("H136").Value = 2121.372
$ws.Range("I136").Value = 1227.3158
$ws.Range("J136").Value = 2829.1667
$ws.Range("K136").Value = 3681.9474
$ws.Range("L136").Value = 8487.500100000001
$ws.Range("M136").Value = -1131.9474
$ws.Range("N136").Value = -13587.5001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449.5
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 499
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 499
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -845
$ws.Range("H62").Value = 47497.5
$ws.Range("J62").Value = 47497.5
$ws.Range("L62").Value = 47497.5
$ws.Range("N62").Value = -48869.5
$ws.Range("H65").Value = 47497.5
$ws.Range("J65").Value = 47497.5
$ws.Range("L65").Value = 142492.5
$ws.Range("N65").Value = -149356.5
$ws.Range("H99").Value = 1950.0416
$ws.Range("J99").Value = 2039.4445
$ws.Range("L99").Value = 2039.4445
$ws.Range("N99").Value = -5035.4445
$ws.Range("H107").Value = 55603628
$ws.Range("J107").Value = 702
$ws.Range("L107").Value = 702
$ws.Range("N107").Value = -4542
$ws.Range("H134").Value = 7369.381
$ws.Range("I134").Value = 6926.9414
$ws.Range("J134").Value = 9249.75
$ws.Range("K134").Value = 20780.8242
$ws.Range("L134").Value = 27749.25
$ws.Range("M134").Value = -18245.8242
$ws.Range("N134").Value = -32819.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2500725
$ws.Range("I4").Value = 2500725
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2500725
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2500613
$ws.Range("H31").Value = 26599.824
$ws.Range("I31").Value = 1229.4814
$ws.Range("J31").Value = 49433.133
$ws.Range("K31").Value = 1229.4814
$ws.Range("L31").Value = 49433.133
$ws.Range("M31").Value = -934.4813999999999
$ws.Range("N31").Value = -50023.133
$ws.Range("H34").Value = 26599.824
$ws.Range("I34").Value = 1229.4814
$ws.Range("J34").Value = 49433.133
$ws.Range("K34").Value = 1229.4814
$ws.Range("L34").Value = 49433.133
$ws.Range("M34").Value = -1027.4814
$ws.Range("N34").Value = -49837.133
$ws.Range("H99").Value = 15460.846
$ws.Range("I99").Value = 6467.75
$ws.Range("J99").Value = 19457.777
$ws.Range("K99").Value = 6467.75
$ws.Range("L99").Value = 19457.777
$ws.Range("M99").Value = -4969.75
$ws.Range("N99").Value = -22453.777
$ws.Range("H107").Value = 789.4091
$ws.Range("I107").Value = 1001.0833
$ws.Range("J107").Value = 535.4
$ws.Range("K107").Value = 1001.0833
$ws.Range("L107").Value = 535.4
$ws.Range("M107").Value = 918.9167
$ws.Range("N107").Value = -4375.4
$ws.Range("H126").Value = 15460.846
$ws.Range("I126").Value = 6467.75
$ws.Range("J126").Value = 19457.777
$ws.Range("K126").Value = 19403.25
$ws.Range("L126").Value = 58373.33099999999
$ws.Range("M126").Value = -16933.25
$ws.Range("N126").Value = -63313.33099999999
$ws.Range("H134").Value = 1330.75
$ws.Range("I134").Value = 1221.7646
$ws.Range("J134").Value = 1499.1818
$ws.Range("K134").Value = 3665.2938
$ws.Range("L134").Value = 4497.5454
$ws.Range("M134").Value = -1130.2938
$ws.Range("N134").Value = -9567.545399999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1424.5
$ws.Range("I4").Value = 333
$ws.Range("J4").Value = 1788.3334
$ws.Range("K4").Value = 999
$ws.Range("L4").Value = 5365.0002
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -5589.0002
$ws.Range("H131").Value = 9602.841
$ws.Range("J131").Value = 9670.529
$ws.Range("L131").Value = 29011.587
$ws.Range("N131").Value = -39091.587

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 14498.5
$ws.Range("J48").Value = 14498.5
$ws.Range("L48").Value = 14498.5
$ws.Range("N48").Value = -15468.5
$ws.Range("H95").Value = 19800
$ws.Range("J95").Value = 19800
$ws.Range("L95").Value = 19800
$ws.Range("N95").Value = -25292
$ws.Range("H97").Value = 66668550
$ws.Range("I97").Value = 83335370
$ws.Range("J97").Value = 1270
$ws.Range("K97").Value = 83335370
$ws.Range("L97").Value = 1270
$ws.Range("M97").Value = -83334874
$ws.Range("N97").Value = -2262
$ws.Range("H126").Value = 3788
$ws.Range("I126").Value = 3776
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 11328
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -8858
$ws.Range("N126").Value = -16340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 82181.03
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 15330
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 15330
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -15554
$ws.Range("H61").Value = 2185.8572
$ws.Range("I61").Value = 2340.6
$ws.Range("J61").Value = 1799
$ws.Range("K61").Value = 2340.6
$ws.Range("L61").Value = 1799
$ws.Range("M61").Value = -2138.6
$ws.Range("N61").Value = -2203
$ws.Range("H113").Value = 2185.8572
$ws.Range("I113").Value = 2340.6
$ws.Range("J113").Value = 1799
$ws.Range("K113").Value = 2340.6
$ws.Range("L113").Value = 1799
$ws.Range("M113").Value = -170.5999999999999
$ws.Range("N113").Value = -6139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 38888
$ws.Range("J31").Value = 38888
$ws.Range("L31").Value = 38888
$ws.Range("N31").Value = -39584
$ws.Range("H56").Value = 35191.332
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 4000
$ws.Range("M56").Value = -3286
$ws.Range("H96").Value = 62501830
$ws.Range("I96").Value = 200002640
$ws.Range("J96").Value = 1463.8182
$ws.Range("K96").Value = 200002640
$ws.Range("L96").Value = 1463.8182
$ws.Range("M96").Value = -200001267
$ws.Range("N96").Value = -4209.8182
